$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "nemad" (ticker symbol) column L previously repeated the full company
# name ("نفت سپاهان"); update it to the actual ticker symbol ("شسپا"),
# which is added as a new shared string. Column M keeps the full company
# name unchanged.
$ws.Range("L2:L45").Value = "شسپا"

# Column L (the symbol column) is narrow text, so it picks up an explicit
# best-fit width once the sheet is re-saved.
$ws.Range("L1").ColumnWidth = 9.166666666666666

# Leave the selection where it was when the file was last saved.
$ws.Range("K4").Select() | Out-Null
